$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New event rows (418-422) — Datum / Event / Location / Stadt / Link
# Values mirror the surrounding rows' pattern: col A is a date serial,
# columns B-E are text (shared strings), col E also carries a hyperlink.
$rows = @(
    @{ Row = 418; Date = 45849; Event = "BASSMANIA FESTIVAL 12 FLOORS 5 CLUBS"; Location = "Favela, Conny Kramer, Fusion, Tryptychon, Sputnik"; Stadt = "Münster"; Link = "https://www.instagram.com/p/DHspQwkNi7E/?igsh=NzA1Yjh6NHVma3M0" },
    @{ Row = 419; Date = 45926; Event = "BASSMANIA PRES. DEEP91 BIRTHDAY WOCHENENDE"; Location = "Favela"; Stadt = "Münster"; Link = "https://www.instagram.com/bassmania_official?igsh=cXViYzk4d3NzMGVr" },
    @{ Row = 420; Date = 45927; Event = "BASSMANIA PRES. DEEP91 BIRTHDAY WOCHENENDE"; Location = "Stollen134"; Stadt = "Dortmund"; Link = "https://www.instagram.com/bassmania_official?igsh=cXViYzk4d3NzMGVr" },
    @{ Row = 421; Date = 45940; Event = "BASSMANIA FESTIVAL 12 FLOORS 5 CLUBS"; Location = "Favela, Conny Kramer, Fusion, Tryptychon, Sputnik"; Stadt = "Münster"; Link = "https://www.instagram.com/p/DH8ogEmNIhM/?igsh=ZHp6Y3R6a3Rtcmh2" },
    @{ Row = 422; Date = 45793; Event = "HARD.NOISE"; Location = "Favela"; Stadt = "Münster"; Link = "https://www.instagram.com/p/DJCDm7kNc05/?igsh=MWJjN3praXp3NTRndQ==" }
)

# A style-3 cell elsewhere in the sheet (text, left border/fill, default font) —
# used below to restore the pre-existing cell style on the link column after
# Hyperlinks.Add overwrites it with its own "visited/unvisited link" style.
$formatDonor = $ws.Range("D417")

foreach ($r in $rows) {
    $rowNum = $r.Row

    $dateCell = $ws.Cells.Item($rowNum, 1)
    $eventCell = $ws.Cells.Item($rowNum, 2)
    $locationCell = $ws.Cells.Item($rowNum, 3)
    $stadtCell = $ws.Cells.Item($rowNum, 4)
    $linkCell = $ws.Cells.Item($rowNum, 5)

    # Date column keeps its existing date-number style; just set the value.
    $dateCell.Value = $r.Date

    # Text columns: mark them as Text-formatted (matches style used by
    # every other populated row) before assigning the string values.
    $eventCell.NumberFormat = "@"
    $eventCell.Value = $r.Event

    $locationCell.NumberFormat = "@"
    $locationCell.Value = $r.Location

    $stadtCell.NumberFormat = "@"
    $stadtCell.Value = $r.Stadt

    $linkCell.NumberFormat = "@"
    $linkCell.Value = $r.Link

    # Register the real hyperlink (adds the relationship + <hyperlink> entry).
    $ws.Hyperlinks.Add($linkCell, $r.Link, "", "", $r.Link) | Out-Null

    # Hyperlinks.Add stamps its own font/style onto the cell; paste back the
    # plain text-cell formatting so the link column matches its neighbours.
    $formatDonor.Copy() | Out-Null
    $linkCell.PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = $false

Write-Host "Added rows 418-422 with BASSMANIA / HARD.NOISE events and their Instagram links"
